# "Correcciones a scripts para datos de prueba."
# Fix the test-data generation sheet: rename the tab to match its real
# purpose and correct a handful of wrong id_usuario values that had been
# left at the default "3" in the generated INSERT statements.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The tab was mislabeled "usuario"; it is really the prog_taller sheet.
$ws.Name = "prog_taller"

# Column G holds id_usuario, used by the generated INSERT statement in
# column H. Several rows still had the placeholder value 3 - correct them.
$ws.Range("G5").Value = 5
$ws.Range("G8").Value = 7
$ws.Range("G11").Value = 5
$ws.Range("G13").Value = 5
$ws.Range("G16").Value = 8
$ws.Range("G19").Value = 10

# Reflect where the author ended up looking while checking the generated
# SQL statements.
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 19
$ws.Range("H3:H20").Select() | Out-Null
